# On Screen Actions - update the "nov" test data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("nov")

# Row 4 used to describe a "type" = Test/Complete split across C4/D4.
# Collapse that into a single "product" = "TestComplete" row, clearing D4.
$ws.Range("B4").Value = "product"
$ws.Range("C4").Value = "TestComplete"
$ws.Range("D4").ClearContents()

# Move the active selection to C5 (was D32).
$ws.Activate()
$ws.Range("C5").Select()
